$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '51.811.88'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = '2.838.36'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '351.31'
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('D6').Value = '113.38'
$ws.Range('E6').Value = '  +5.01%  '
$ws.Range('E7').Value = '  +2.01%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.622'
$ws.Range('E9').Value = '  +6.02%  '
$ws.Range('D10').Value = '40.28'
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = '0.0852'
$ws.Range('E12').Value = '  +2.32%  '
$ws.Range('D13').Value = '20.03'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').Value = '7.81'
$ws.Range('E14').Value = '  +3.54%  '
$ws.Range('D15').Value = '3.271.73'
$ws.Range('E15').Value = '  +1.95%  '
$ws.Range('D16').Value = '0.976'
$ws.Range('E16').Value = '  +5.58%  '
$ws.Range('D17').Value = '2.829.10'
$ws.Range('E17').Value = '  +1.97%  '
$ws.Range('D18').Value = '51.856.49'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '3.46'
$ws.Range('E19').Value = '  +11.92%  '
$ws.Range('D20').Value = '7.63'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '13.38'
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').Value = '0.0₃0975'
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').Value = '70.61'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('D24').Value = '269.34'
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').Value = '2.77'
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('D26').Value = '26.36'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '39.14'
$ws.Range('E29').Value = '  +7.01%  '
$ws.Range('D30').Value = '10.62'
$ws.Range('E30').Value = '  +3.92%  '
$ws.Range('D31').Value = '2.26'
$ws.Range('E31').Value = '  +2.44%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '6.27'
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = '52.84'
$ws.Range('E33').Value = '  +2.02%  '
$ws.Range('D34').Value = '0.0457'
$ws.Range('E34').Value = '  +2.81%  '
$ws.Range('D35').Value = '0.0902'
$ws.Range('E35').Value = '  +9.36%  '
$ws.Range('D36').Value = '5.66'
$ws.Range('E36').Value = '  +2.40%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').Value = '19.12'
$ws.Range('E38').Value = '  +4.72%  '
$ws.Range('D39').Value = '3.24'
$ws.Range('E39').Value = '  +3.17%  '
$ws.Range('D40').Value = '2.02'
$ws.Range('E40').Value = '  +3.10%  '
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('D43').Value = '121.83'
$ws.Range('E43').Value = '  +0.80%  '
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('D45').Value = '22.11'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('D46').Value = '2.184.60'
$ws.Range('E46').Value = '  +4.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.50'
$ws.Range('E47').Value = '  +7.83%  '
$ws.Range('E48').Value = '  +7.34%  '
$ws.Range('E49').Value = '  +26.53%  '
$ws.Range('D50').Value = '0.975'
$ws.Range('E50').Value = '  +7.90%  '
$ws.Range('D51').Value = '5.55'
$ws.Range('E51').Value = '  +2.30%  '
